$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename cluster label in A3 from "MuSCs" to "Inflammatory-Mac"
$ws.Range("A3").Value = "Inflammatory-Mac"

# Row 2 (Hcrt -> Hcrtr2, ECs sending cluster) updated TPM-derived values
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.39064
$ws.Range("H2").Value = 1.17192
$ws.Range("I2").Value = 0.5296693860025763
$ws.Range("J2").Value = 0.5296693860025762
$ws.Range("Q2").Value = 0.002938394080000001
$ws.Range("R2").Value = 0.02644554672
$ws.Range("S2").Value = 0.5296693860025763
$ws.Range("T2").Value = 0.5296693860025762

# Row 3 (Inflammatory-Mac sending cluster) updated TPM-derived values
$ws.Range("G3").Value = 0.1178916666666667
$ws.Range("H3").Value = 0.353675
$ws.Range("I3").Value = 0.1598494949266683
$ws.Range("J3").Value = 0.1598494949266683
$ws.Range("Q3").Value = 0.0008867811166666667
$ws.Range("R3").Value = 0.007981030050000001
$ws.Range("S3").Value = 0.1598494949266683
$ws.Range("T3").Value = 0.1598494949266683

# Row 4 (Resolving-Mac sending cluster) updated TPM-derived values
$ws.Range("G4").Value = 0.228985
$ws.Range("H4").Value = 0.686955
$ws.Range("I4").Value = 0.3104811190707554
$ws.Range("J4").Value = 0.3104811190707554
$ws.Range("Q4").Value = 0.00172242517
$ws.Range("R4").Value = 0.01550182653
$ws.Range("S4").Value = 0.3104811190707554
$ws.Range("T4").Value = 0.3104811190707554
